# Adds a "time_taken" column (F) with per-row timestamps to the "data" sheet,
# matching the author's commit: "Updated panelApp panels to contain
# time_taken for metadata purposes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same bold/centered/bordered style as the other headers (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps (as plain text, matching column E's un-styled data cells)
$timestamps = @{
    2  = "2021-10-05 13:40:44.991362"
    3  = "2021-10-05 13:40:44.991373"
    4  = "2021-10-05 13:40:44.991377"
    5  = "2021-10-05 13:40:44.991380"
    6  = "2021-10-05 13:40:44.991384"
    7  = "2021-10-05 13:40:44.991387"
    8  = "2021-10-05 13:40:44.991390"
    9  = "2021-10-05 13:40:44.991393"
    10 = "2021-10-05 13:40:44.991396"
    11 = "2021-10-05 13:40:44.991399"
    12 = "2021-10-05 13:40:44.991402"
    13 = "2021-10-05 13:40:44.991405"
    14 = "2021-10-05 13:40:44.991408"
    15 = "2021-10-05 13:40:44.991411"
    16 = "2021-10-05 13:40:44.991413"
    17 = "2021-10-05 13:40:44.991417"
    18 = "2021-10-05 13:40:44.991420"
    19 = "2021-10-05 13:40:44.991423"
    20 = "2021-10-05 13:40:44.991426"
    21 = "2021-10-05 13:40:44.991429"
    22 = "2021-10-05 13:40:44.991432"
    23 = "2021-10-05 13:40:44.991435"
    24 = "2021-10-05 13:40:44.991438"
    25 = "2021-10-05 13:40:44.991440"
    26 = "2021-10-05 13:40:44.991444"
    27 = "2021-10-05 13:40:44.991447"
    28 = "2021-10-05 13:40:44.991450"
    29 = "2021-10-05 13:40:44.991453"
    30 = "2021-10-05 13:40:44.991456"
    31 = "2021-10-05 13:40:44.991459"
    32 = "2021-10-05 13:40:44.991462"
    33 = "2021-10-05 13:40:44.991465"
    34 = "2021-10-05 13:40:44.991468"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 6).Value = $timestamps[$row]
}
